# Generate Report for Handoff
# Adds a new row (file cd7e7da9-...-eed8....md, now "Ready for handoff") to the
# Overview sheet and the two per-locale sheets (zh-cn, de-de), mirroring the
# existing row for the 68f9f1e9-... file already present in row 2 of each
# table.

$wb = $excel.ActiveWorkbook

# ---- shared literal strings (kept identical across sheets so the shared
#      string table de-dupes exactly like the source workbook) ----
$S_MD_EXT            = ".md"
$S_EMPTY             = ""
$S_NEW_SRC_MD        = "cd7e7da9-44ac-4528-9e66-10281258eed8ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$S_NEW_DISPLAY_MD    = "e2e\cd7e7da9-44ac-4528-9e66-10281258eed8ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$S_READY_FOR_HANDOFF = "Ready for handoff"
$S_DATE_1229_54      = "2016-08-12 12:29:54"
$S_E2E               = "e2e"
$S_HT                = "ht"
$S_FALSE             = "False"
$S_EPOCH             = "0001-01-01 00:00:00"
$S_TRUE              = "True"
$S_NEW_ZHCN_XLF      = "cd7e7da9-44ac-4528-9e66-10281258eed8ooooooooooooooooooooooooooooooooooooooooooo.c7184a0924ac8b062d7f5645ca2104b2496493a0.zh-cn.xlf"
$S_DATE_1229_48      = "2016-08-12 12:29:48"
$S_NEW_DEDE_XLF      = "cd7e7da9-44ac-4528-9e66-10281258eed8ooooooooooooooooooooooooooooooooooooooooooo.c7184a0924ac8b062d7f5645ca2104b2496493a0.de-de.xlf"

$NEW_URL = "https://github.com/OpenLocalizationTestOrg/oltest/blob/ad4215389dc0a63159c7eb665efe705f50a08b85/e2e/cd7e7da9-44ac-4528-9e66-10281258eed8ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# =====================================================================
# Sheet 1: Overview  -> new row 3
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $S_NEW_SRC_MD
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $NEW_URL, "", "", $S_NEW_DISPLAY_MD)
$wsOverview.Cells.Item(3, 3).Value = $S_MD_EXT
$wsOverview.Cells.Item(3, 4).Value = $S_EMPTY
$wsOverview.Cells.Item(3, 5).Value = $S_READY_FOR_HANDOFF
$wsOverview.Cells.Item(3, 6).Value = $S_READY_FOR_HANDOFF
$wsOverview.Cells.Item(3, 7).Value = $S_DATE_1229_54

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# =====================================================================
# Sheet 2: zh-cn -> new row 3
# =====================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $NEW_URL, "", "", $S_NEW_SRC_MD)
$wsZhCn.Cells.Item(3, 2).Value  = $S_MD_EXT
$wsZhCn.Cells.Item(3, 3).Value  = $S_READY_FOR_HANDOFF
$wsZhCn.Cells.Item(3, 4).Value  = $S_E2E
$wsZhCn.Cells.Item(3, 5).Value  = $S_HT
$wsZhCn.Cells.Item(3, 6).Value  = $S_FALSE
$wsZhCn.Cells.Item(3, 7).Value  = $S_NEW_ZHCN_XLF
$wsZhCn.Cells.Item(3, 8).Value  = $S_DATE_1229_48
$wsZhCn.Cells.Item(3, 9).Value  = $S_EMPTY
$wsZhCn.Cells.Item(3, 10).Value = $S_EMPTY
$wsZhCn.Cells.Item(3, 11).Value = $S_EPOCH
$wsZhCn.Cells.Item(3, 12).Value = $S_EMPTY
$wsZhCn.Cells.Item(3, 13).Value = $S_TRUE
$wsZhCn.Cells.Item(3, 14).Value = $S_EMPTY
$wsZhCn.Cells.Item(3, 15).Value = $S_FALSE
$wsZhCn.Cells.Item(3, 16).Value = $S_EMPTY

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# =====================================================================
# Sheet 3: de-de -> new row 3
# =====================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $NEW_URL, "", "", $S_NEW_SRC_MD)
$wsDeDe.Cells.Item(3, 2).Value  = $S_MD_EXT
$wsDeDe.Cells.Item(3, 3).Value  = $S_READY_FOR_HANDOFF
$wsDeDe.Cells.Item(3, 4).Value  = $S_E2E
$wsDeDe.Cells.Item(3, 5).Value  = $S_HT
$wsDeDe.Cells.Item(3, 6).Value  = $S_FALSE
$wsDeDe.Cells.Item(3, 7).Value  = $S_NEW_DEDE_XLF
$wsDeDe.Cells.Item(3, 8).Value  = $S_DATE_1229_54
$wsDeDe.Cells.Item(3, 9).Value  = $S_EMPTY
$wsDeDe.Cells.Item(3, 10).Value = $S_EMPTY
$wsDeDe.Cells.Item(3, 11).Value = $S_EPOCH
$wsDeDe.Cells.Item(3, 12).Value = $S_EMPTY
$wsDeDe.Cells.Item(3, 13).Value = $S_TRUE
$wsDeDe.Cells.Item(3, 14).Value = $S_EMPTY
$wsDeDe.Cells.Item(3, 15).Value = $S_FALSE
$wsDeDe.Cells.Item(3, 16).Value = $S_EMPTY

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

# Apply the datetime display format to the "date-looking" text cells, matching
# style index 2 used by the existing rows (numFmt "yyyy-mm-dd HH:mm:ss").
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Output "Report row for handoff added."
